# Apply cryptos.xlsx price/volume update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices formatted as text (e.g. "48.077.52"); force text
# storage so Excel does not reinterpret these numeric-looking strings as
# numbers, matching the inlineStr cells produced by the source data feed.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "48.077.52"
$ws.Range("E2").Value = "  +0.17%  "

$ws.Range("D3").Value = "2.502.68"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "320.49"
$ws.Range("E5").Value = "  -0.96%  "

$ws.Range("D6").Value = "107.46"
$ws.Range("E6").Value = "  -2.21%  "

$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").Value = "0.538"
$ws.Range("E9").Value = "  -3.15%  "

$ws.Range("D10").Value = "39.53"
$ws.Range("E10").Value = "  -3.33%  "

$ws.Range("E11").Value = "  +7.61%  "

$ws.Range("D12").Value = "0.0812"
$ws.Range("E12").Value = "  -0.58%  "

$ws.Range("E13").Value = "  -0.10%  "

$ws.Range("D14").Value = "7.10"
$ws.Range("E14").Value = "  -2.24%  "

$ws.Range("D15").Value = "2.893.54"
$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("D16").Value = "2.498.04"
$ws.Range("E16").Value = "  -0.50%  "

$ws.Range("D17").Value = "0.838"
$ws.Range("E17").Value = "  -2.35%  "

$ws.Range("D18").Value = "47.914.46"
$ws.Range("E18").Value = "  +0.01%  "

$ws.Range("D19").Value = "12.94"
$ws.Range("E19").Value = "  -3.19%  "

$ws.Range("D20").Value = "6.73"
$ws.Range("E20").Value = "  +0.92%  "

$ws.Range("E21").Value = "  -0.84%  "

$ws.Range("D22").Value = "2.76"
$ws.Range("E22").Value = "  -1.39%  "

$ws.Range("D23").Value = "277.83"
$ws.Range("E23").Value = "  +11.97%  "

$ws.Range("D24").Value = "71.52"
$ws.Range("E24").Value = "  +0.80%  "

$ws.Range("D25").Value = "2.53"
$ws.Range("E25").Value = "  -1.53%  "

$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("D27").Value = "25.63"
$ws.Range("E27").Value = "  -1.47%  "

$ws.Range("D28").Value = "9.73"
$ws.Range("E28").Value = "  -3.28%  "

$ws.Range("E29").Value = "  -0.89%  "

$ws.Range("D30").Value = "35.10"
$ws.Range("E30").Value = "  -0.56%  "

$ws.Range("D31").Value = "2.10"
$ws.Range("E31").Value = "  -9.30%  "

$ws.Range("D32").Value = "49.48"
$ws.Range("E32").Value = "  -0.60%  "

$ws.Range("D33").Value = "19.51"
$ws.Range("E33").Value = "  -3.30%  "

$ws.Range("E34").Value = "  -0.14%  "

$ws.Range("D35").Value = "5.30"
$ws.Range("E35").Value = "  -1.49%  "

$ws.Range("D36").Value = "0.0781"
$ws.Range("E36").Value = "  -1.26%  "

$ws.Range("E37").Value = "  -1.89%  "

$ws.Range("D38").Value = "4.62"
$ws.Range("E38").Value = "  -1.41%  "

$ws.Range("D39").Value = "2.89"
$ws.Range("E39").Value = "  -3.47%  "

$ws.Range("E40").Value = "  -0.87%  "

$ws.Range("D41").Value = "120.66"
$ws.Range("E41").Value = "  +0.90%  "

$ws.Range("D43").Value = "21.18"
$ws.Range("E43").Value = "  -6.38%  "

$ws.Range("E44").Value = "  +0.24%  "

$ws.Range("D45").Value = "2.007.95"
$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("D46").Value = "3.15"
$ws.Range("E46").Value = "  +2.12%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "1.85"
$ws.Range("E47").Value = "  +0.86%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "2.00"
$ws.Range("E48").Value = "  -1.69%  "

$ws.Range("E49").Value = "  -0.99%  "

$ws.Range("D50").Value = "5.16"
$ws.Range("E50").Value = "  -1.14%  "

$ws.Range("D51").Value = "80.25"
$ws.Range("E51").Value = "  +2.75%  "

# Restore the default (unstyled) cell style on column D now that the
# text values are stored, so formatting matches the original sheet.
$ws.Range("D2:D51").Style = "Normal"
